$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values that are plain numeric-looking strings must be forced to
# text (NumberFormat "@") before assignment, otherwise Excel COM auto-
# coerces them into floating point numbers and mangles formatting like
# trailing zeros (e.g. "1.00" -> 1).

$ws.Cells.Item(2, 4).Value = '69.615.40'
$ws.Cells.Item(2, 5).Value = '  +2.21%  '

$ws.Cells.Item(3, 4).Value = '3.913.52'
$ws.Cells.Item(3, 5).Value = '  +0.29%  '

$ws.Cells.Item(4, 5).Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '526.79'
$ws.Cells.Item(5, 5).Value = '  +8.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '144.49'
$ws.Cells.Item(6, 5).Value = '  -0.83%  '

$ws.Cells.Item(7, 5).Value = '  -1.18%  '

$ws.Cells.Item(8, 5).Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.725'
$ws.Cells.Item(9, 5).Value = '  -1.33%  '

$ws.Cells.Item(10, 5).Value = '  +0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0000336'
$ws.Cells.Item(11, 5).Value = '  -4.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '42.31'
$ws.Cells.Item(12, 5).Value = '  -1.49%  '

$ws.Cells.Item(13, 4).Value = '4.533.51'
$ws.Cells.Item(13, 5).Value = '  +0.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '10.29'
$ws.Cells.Item(14, 5).Value = '  -3.20%  '

$ws.Cells.Item(15, 4).Value = '3.921.55'
$ws.Cells.Item(15, 5).Value = '  -0.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.20'
$ws.Cells.Item(16, 5).Value = '  -0.55%  '

$ws.Cells.Item(17, 5).Value = '  +7.67%  '

$ws.Cells.Item(18, 5).Value = '  -0.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '19.77'
$ws.Cells.Item(19, 5).Value = '  -2.09%  '

$ws.Cells.Item(20, 4).Value = '69.505.32'
$ws.Cells.Item(20, 5).Value = '  +1.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '430.31'
$ws.Cells.Item(21, 5).Value = '  +0.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '3.37'
$ws.Cells.Item(22, 5).Value = '  -5.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '14.27'
$ws.Cells.Item(23, 5).Value = '  -4.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '88.58'
$ws.Cells.Item(24, 5).Value = '  -0.38%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '4.03'
$ws.Cells.Item(25, 5).Value = '  +8.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '11.40'
$ws.Cells.Item(26, 5).Value = '  -0.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.63'
$ws.Cells.Item(27, 5).Value = '  -3.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '36.45'
$ws.Cells.Item(28, 5).Value = '  -2.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '695.08'
$ws.Cells.Item(29, 5).Value = '  -3.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '13.22'
$ws.Cells.Item(30, 5).Value = '  -3.57%  '

$ws.Cells.Item(31, 5).Value = '  -3.44%  '

$ws.Cells.Item(32, 5).Value = '  -3.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '67.06'
$ws.Cells.Item(33, 5).Value = '  +10.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.442'
$ws.Cells.Item(34, 5).Value = '  +11.72%  '

$ws.Cells.Item(35, 5).Value = '  -2.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '40.17'
$ws.Cells.Item(36, 5).Value = '  -2.96%  '

$ws.Cells.Item(37, 4).Value = '0.0₃0852'
$ws.Cells.Item(37, 5).Value = '  -4.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.151'
$ws.Cells.Item(38, 5).Value = '  +3.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.999'
$ws.Cells.Item(39, 5).Value = '  +0.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 5).Value = '  -0.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.0483'
$ws.Cells.Item(41, 5).Value = '  -2.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.10'
$ws.Cells.Item(42, 5).Value = '  +4.01%  '

$ws.Cells.Item(43, 2).Value = 'ThetaToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D43").NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.09'
$ws.Cells.Item(43, 5).Value = '  -0.10%  '

$ws.Cells.Item(44, 2).Value = 'Fetch.AI'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.80'
$ws.Cells.Item(44, 5).Value = '  -7.77%  '

$ws.Cells.Item(45, 5).Value = '  -1.04%  '

$ws.Cells.Item(46, 5).Value = '  -0.39%  '

$ws.Cells.Item(47, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(47, 4).Value = '0.0₆0354'
$ws.Cells.Item(47, 5).Value = '  +7.44%  '

$ws.Cells.Item(48, 2).Value = 'Stacks'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '3.01'
$ws.Cells.Item(48, 5).Value = '  +7.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '26.70'
$ws.Cells.Item(49, 5).Value = '  +5.81%  '

$ws.Cells.Item(50, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '3.29'
$ws.Cells.Item(50, 5).Value = '  -3.34%  '

$ws.Cells.Item(51, 2).Value = 'Maker'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(51, 4).Value = '2.727.79'
$ws.Cells.Item(51, 5).Value = '  +11.56%  '
